# Deploying to gh-pages: add the 2021 data column (R) to the sheet and
# adjust header row heights / selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights for the title rows -----------------------------------
$ws.Rows.Item(1).RowHeight = 41.25
$ws.Rows.Item(2).RowHeight = 15

# --- Copy the Q column (2020) formatting into the new R column (2021) -
# This brings along the per-cell styles (number format, borders, etc.)
# for every row that has a Q entry, rows 3-38.
$ws.Range("Q3:Q38").Copy()
$null = $ws.Range("R3:R38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New column header (year 2021) ------------------------------------
$ws.Range("R4").Value = 2021

# --- New 2021 data values ----------------------------------------------
$ws.Range("R6").Value = 88.796593100633856
$ws.Range("R7").Value = 86.908583391486388
$ws.Range("R8").Value = 89.680106631122953
$ws.Range("R9").Value = 95.775910364145659
$ws.Range("R10").Value = 96.517042279754136
$ws.Range("R11").Value = 90.311530128242666
$ws.Range("R12").Value = 90.746324915190343
$ws.Range("R13").Value = 90.894107952204379
$ws.Range("R14").Value = 81.065680730752504
$ws.Range("R15").Value = 85.088888888888889

$ws.Range("R17").Value = 93.37839883628321
$ws.Range("R18").Value = 93.091416608513612
$ws.Range("R19").Value = 94.815061646117954
$ws.Range("R20").Value = 100.53781512605042
$ws.Range("R21").Value = 100.33525796237662
$ws.Range("R22").Value = 93.78989283832054
$ws.Range("R23").Value = 95.401432340746325
$ws.Range("R24").Value = 92.308748798242007
$ws.Range("R25").Value = 89.338842975206617
$ws.Range("R26").Value = 87.955555555555549

$ws.Range("R28").Value = 89.631204460036727
$ws.Range("R29").Value = 89.204466154919743
$ws.Range("R30").Value = 84.751749416861045
$ws.Range("R31").Value = 96.201680672268907
$ws.Range("R32").Value = 95.567144719687093
$ws.Range("R33").Value = 91.330444457457389
$ws.Range("R34").Value = 91.368262344515642
$ws.Range("R35").Value = 92.345373803964662
$ws.Range("R36").Value = 88.660287081339717
$ws.Range("R37").Value = 84.944444444444443

# Rows 5, 16, 27 and 38 stay blank in the R column (same as their Q
# counterparts), only formatting was copied for those.

# --- Selection moves to R3, matching the author's last-active cell ----
$null = $ws.Range("R3").Select()
